$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.0122
$ws.Range("E2").Value = -0.003500000000000003
$ws.Range("G2").Value = 0.07608255766895992
$ws.Range("H2").Value = 0.07608255766895992
$ws.Range("I2").Value = 0.06304298696042929
$ws.Range("J2").Value = 0.04840272087365585
$ws.Range("K2").Value = 516.7
$ws.Range("L2").Value = 0.02715657470817325
$ws.Range("M2").Value = 62.1
$ws.Range("N2").Value = 0.01097173144876325
$ws.Range("O2").Value = 0.1201857944648732
$ws.Range("P2").Value = 62.1
$ws.Range("Q2").Value = 0.01097173144876325
$ws.Range("R2").Value = 0.1201857944648732
$ws.Range("U2").Value = 2485.2
$ws.Range("V2").Value = 0.4390812720848056
$ws.Range("W2").Value = 0.05531379304623442
$ws.Range("X2").Value = 0.05550309912554552
$ws.Range("Y2").Value = -0.0001893060793110996
$ws.Range("Z2").Value = 1.942828260137033
$ws.Range("AA2").Value = 0.08746185060957173
$ws.Range("AB2").Value = 0.0428053462132216
$ws.Range("AC2").Value = 0.04465650439635014
$ws.Range("AD2").Value = 2818.2
$ws.Range("AE2").Value = 0.0
$ws.Range("AF2").Value = 2818.2
$ws.Range("AG2").Value = 333.0
$ws.Range("AH2").Value = 0.3324054634238399
$ws.Range("AI2").Value = 0.2258191171403617
$ws.Range("AJ2").Value = 0.05556482562990155
$ws.Range("AK2").Value = 0.03331765835893023
$ws.Range("AL2").Value = 190.9
$ws.Range("AM2").Value = 190.9
$ws.Range("AN2").Value = 2.026461494211548
$ws.Range("AO2").Value = 6.283394447354635
$ws.Range("AP2").Value = 0.2394477601208025
$ws.Range("AQ2").Value = 6.283394447354635

# Row 3
$ws.Range("D3").Value = 0.0273
$ws.Range("E3").Value = 0.117
$ws.Range("G3").Value = 0.07995483688858422
$ws.Range("H3").Value = 0.07995483688858422
$ws.Range("I3").Value = 0.07278519295132868
$ws.Range("J3").Value = 0.0532001803288512
$ws.Range("K3").Value = 328.5
$ws.Range("L3").Value = 0.02649300375015122
$ws.Range("U3").Value = 1849.5
$ws.Range("V3").Value = 0.5678887251289609
$ws.Range("W3").Value = 0.05933994472443505
$ws.Range("X3").Value = 0.05638405988183222
$ws.Range("Y3").Value = 0.002955884842602834
$ws.Range("Z3").Value = 2.221455828869341
$ws.Range("AA3").Value = 0.1181818506884266
$ws.Range("AB3").Value = 0.04279434891613403
$ws.Range("AC3").Value = 0.07538750177229255
$ws.Range("AD3").Value = 1717.9
$ws.Range("AE3").Value = 0.0
$ws.Range("AF3").Value = 1717.9
$ws.Range("AG3").Value = -131.5999999999999
$ws.Range("AH3").Value = 0.3453273564235029
$ws.Range("AI3").Value = 0.2243480077834224
$ws.Range("AJ3").Value = -0.04210930500447968
$ws.Range("AK3").Value = -0.02265918247873548
$ws.Range("AL3").Value = 127.5
$ws.Range("AM3").Value = 127.5
$ws.Range("AN3").Value = 1.695351820783578
$ws.Range("AO3").Value = 7.078431372549019
$ws.Range("AP3").Value = -0.1298726931806967
$ws.Range("AQ3").Value = 7.078431372549019

# Row 4
$ws.Range("D4").Value = -0.0517
$ws.Range("E4").Value = -0.124
$ws.Range("G4").Value = 0.06883751810719459
$ws.Range("H4").Value = 0.06883751810719459
$ws.Range("I4").Value = 0.04481530661516176
$ws.Range("J4").Value = 0.03605956930455808
$ws.Range("K4").Value = 188.2
$ws.Range("L4").Value = 0.02839811685176243
$ws.Range("M4").Value = 62.1
$ws.Range("N4").Value = 0.02584054593874834
$ws.Range("O4").Value = 0.3299681190223167
$ws.Range("P4").Value = 62.1
$ws.Range("Q4").Value = 0.02584054593874834
$ws.Range("R4").Value = 0.3299681190223167
$ws.Range("U4").Value = 635.7
$ws.Range("V4").Value = 0.2645223035952065
$ws.Range("W4").Value = 0.05128764136803379
$ws.Range("X4").Value = 0.05462213836925882
$ws.Range("Y4").Value = -0.003334497001225033
$ws.Range("Z4").Value = 1.573558742520657
$ws.Range("AA4").Value = 0.0567418505307169
$ws.Range("AB4").Value = 0.04281634351030917
$ws.Range("AC4").Value = 0.01392550702040773
$ws.Range("AD4").Value = 1100.3
$ws.Range("AE4").Value = 0.0
$ws.Range("AF4").Value = 1100.3
$ws.Range("AG4").Value = 464.5999999999999
$ws.Range("AH4").Value = 0.3140573712002283
$ws.Range("AI4").Value = 0.2281549371708207
$ws.Range("AJ4").Value = 0.1620057186693633
$ws.Range("AK4").Value = 0.1109651532159832
$ws.Range("AL4").Value = 63.4
$ws.Range("AM4").Value = 63.4
$ws.Range("AN4").Value = 2.915474297827239
$ws.Range("AO4").Value = 4.684542586750789
$ws.Range("AP4").Value = 1.23105458399576
$ws.Range("AQ4").Value = 4.684542586750789
